# Apply cryptos list update (prices/volume %) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '67.876.43'
$cell.Style = 'Normal'
$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  -1.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.825.54'
$cell.Style = 'Normal'
$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  -2.14%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.03%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '600.15'
$cell.Style = 'Normal'

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '169.77'
$cell.Style = 'Normal'
$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  -0.10%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '3.822.34'
$cell.Style = 'Normal'
$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -2.21%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.531'
$cell.Style = 'Normal'
$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -0.13%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.166'
$cell.Style = 'Normal'
$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  -0.73%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  +1.17%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  +0.20%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.0000276'
$cell.Style = 'Normal'
$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +8.34%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '37.11'
$cell.Style = 'Normal'
$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -0.59%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '4.463.81'
$cell.Style = 'Normal'
$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -2.11%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '3.840.73'
$cell.Style = 'Normal'
$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  -1.71%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '18.64'
$cell.Style = 'Normal'
$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  +2.30%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '67.953.54'
$cell.Style = 'Normal'
$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  -0.84%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +0.38%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  +0.20%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '10.89'
$cell.Style = 'Normal'
$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  +0.43%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '468.99'
$cell.Style = 'Normal'
$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  -1.26%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '0.742'
$cell.Style = 'Normal'
$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -8.79%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '83.74'
$cell.Style = 'Normal'
$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -0.19%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '2.32'
$cell.Style = 'Normal'
$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  +1.97%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '10.42'
$cell.Style = 'Normal'
$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  +3.72%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  -0.13%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -1.63%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '3.970.65'
$cell.Style = 'Normal'
$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -2.05%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '7.78'
$cell.Style = 'Normal'
$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -1.79%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '2.29'
$cell.Style = 'Normal'
$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  -1.74%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '30.82'
$cell.Style = 'Normal'
$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  -2.51%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '9.38'
$cell.Style = 'Normal'
$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -1.36%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '3.788.50'
$cell.Style = 'Normal'
$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -2.26%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '3.89'
$cell.Style = 'Normal'
$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +5.04%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.106'
$cell.Style = 'Normal'
$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +1.04%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '6.01'
$cell.Style = 'Normal'
$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +0.90%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -1.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.02'
$cell.Style = 'Normal'

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +1.51%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B44')
$cell.NumberFormat = '@'
$cell.Value = 'Cosmos'
$cell.Style = 'Normal'
$cell = $ws.Range('C44')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '8.82'
$cell.Style = 'Normal'
$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  +1.44%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B45')
$cell.NumberFormat = '@'
$cell.Value = 'USDe'
$cell.Style = 'Normal'
$cell = $ws.Range('C45')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  -0.02%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '1.98'
$cell.Style = 'Normal'
$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  -1.24%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '411.84'
$cell.Style = 'Normal'
$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -4.15%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  -1.22%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.000286'
$cell.Style = 'Normal'
$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -4.95%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '143.06'
$cell.Style = 'Normal'
$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  -0.79%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -0.20%  '
$cell.Style = 'Normal'
